$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.830.60"
$ws.Range("E2").Value = "  +5.10%  "
$ws.Range("D3").Value = "3.331.09"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("E4").Value = "  +0.10%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "409.47"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.52%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "111.69"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.25%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.582"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +4.35%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +0.97%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "39.75"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.80%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0986"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.51%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.144"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("D13").Value = "3.841.03"
$ws.Range("E13").Value = "  +1.87%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "8.48"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +4.54%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "19.21"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "3.345.57"
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").Value = "59.588.78"
$ws.Range("E18").Value = "  +4.92%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "10.62"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.39%  "
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("E21").Value = "  +2.54%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "13.07"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.03%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "303.53"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.40%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "75.39"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.46%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.18"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.183"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +8.32%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "28.50"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "4.46"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.68%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.77"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +29.95%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.82"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.32%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.42"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("E32").Value = "  +3.29%  "
$ws.Range("E33").Value = "  +0.16%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "11.52"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +4.50%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "39.54"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +4.51%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.0506"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +4.36%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "51.76"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.37%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.15"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.89%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.995"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.32%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.39"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -4.94%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "139.04"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +3.73%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.123"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.38%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.92"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("E45").Value = "  -2.16%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "16.79"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -4.23%  "
$ws.Range("E47").Value = "  +8.52%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "22.43"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("D49").Value = "2.197.07"
$ws.Range("E49").Value = "  +2.31%  "
$ws.Range("E50").Value = "  -0.05%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.02"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.76%  "
